# Update "JASOTAKO DIRUA" (funding) figures for several days and move the
# active selection, matching the author's manual data-entry session.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja 1")

# Abenduak 12 (row 8) - LH column now has a value
$ws.Range("D8").Value = 321

# Abenduak 13 (row 9) - LH column now has a value
$ws.Range("D9").Value = 454.1

# Abenduak 16 (row 12) - LH column now has a value
$ws.Range("D12").Value = 340

# Abenduak 17 (row 13) - all three columns now have values
$ws.Range("C13").Value = 25
$ws.Range("D13").Value = 322
$ws.Range("E13").Value = 125

# Leave the selection where the author left it when they saved the file
[void]$ws.Range("E11").Select()
